# Re-order the per-player data (columns B:K) for three groups of rows on
# the roster sheet, while leaving column A (the running index) untouched.
# Mapping of "new row <- old row" content:
#   row 5 <- row 6, row 6 <- row 7, row 7 <- row 5   (3-cycle)
#   row 9 <-> row 10                                  (swap)
#   row 12 <-> row 13                                 (swap)
#
# Columns C,D,E,G,H,I,J hold text (player name / position / height /
# birth date / nationality / experience / college) even when the text
# looks numeric (e.g. Exp = "9"); a plain `.Value = "9"` would be
# auto-coerced to a real number by Excel's smart-entry parsing, which
# would change the cell's stored type (and drop it out of the shared
# string table) relative to the source workbook. Forcing the cell to a
# text number-format before the write -- then clearing formatting again
# so no stray style index is left behind -- keeps the value textual
# while leaving the cell's appearance/style exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCols = @("C", "D", "E", "G", "H", "I", "J")
$numCols  = @("B", "F")
$urlCols  = @("K")
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

function Get-RowValues($ws, $row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value2
    }
    return $vals
}

function Set-RowValues($ws, $row, $textCols, $numCols, $urlCols, $vals) {
    foreach ($c in $textCols) {
        $v = $vals[$c]
        $cell = $ws.Range("$c$row")
        if ($v -eq $null) {
            $cell.ClearContents()
        } else {
            $cell.NumberFormat = "@"
            $cell.Value = [string]$v
            $cell.ClearFormats()
        }
    }
    foreach ($c in ($numCols + $urlCols)) {
        $v = $vals[$c]
        $cell = $ws.Range("$c$row")
        if ($v -eq $null) {
            $cell.ClearContents()
        } else {
            $cell.Value = $v
        }
    }
}

# Snapshot the "before" values for every row that will move.
$row5 = Get-RowValues $ws 5 $cols
$row6 = Get-RowValues $ws 6 $cols
$row7 = Get-RowValues $ws 7 $cols
$row9 = Get-RowValues $ws 9 $cols
$row10 = Get-RowValues $ws 10 $cols
$row12 = Get-RowValues $ws 12 $cols
$row13 = Get-RowValues $ws 13 $cols

# 3-cycle across rows 5, 6, 7.
Set-RowValues $ws 5 $textCols $numCols $urlCols $row6
Set-RowValues $ws 6 $textCols $numCols $urlCols $row7
Set-RowValues $ws 7 $textCols $numCols $urlCols $row5

# Swap rows 9 and 10.
Set-RowValues $ws 9 $textCols $numCols $urlCols $row10
Set-RowValues $ws 10 $textCols $numCols $urlCols $row9

# Swap rows 12 and 13.
Set-RowValues $ws 12 $textCols $numCols $urlCols $row13
Set-RowValues $ws 13 $textCols $numCols $urlCols $row12
